$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: add a thin bottom border under the whole row (new separator
#     border style), turning B3/C3/D3/E3 (and the previously-empty A3)
#     into bordered cells.
$ws.Range("A3:E3").Borders.Item(9).LineStyle = 1
$ws.Range("A3:E3").Borders.Item(9).Weight = 2

# --- Rows 4-5: new script block (filename + eng + translated + converted).
#     Values are written in the same column-interleaved order the original
#     authoring tool used for rows 2-3 (C, A, C, D, D, E, E) so the shared
#     string table is appended in the matching sequence.
$ws.Range("C4").Value = " If we make a new discovery,\nwe\'ll make sure to let everyone know."
$ws.Range("A4").Value = "SCRIPT/P01P04A/um0725.ssb"
$ws.Range("C5").Value = " Please look forward to it!"
$ws.Range("D4").Value = " Если мы совершим новое открытие,\nмы всем об этом расскажем."
$ws.Range("D5").Value = " Следите за новыми новостями!"
$ws.Range("E4").Value = " Åòìé íú òïâåñšéí îïâïå ïóëñúóéå,\níú âòåí ïá üóïí ñàòòëàçåí."
$ws.Range("E5").Value = " Òìåäéóå èà îïâúíé îïâïòóÿíé!"

$ws.Range("B4").Value = 98
$ws.Range("B5").Value = 101
$ws.Rows.Item(4).RowHeight = 43.2

# --- Selection moves to B3 (matches the author's saved cursor position)
$ws.Range("B3").Select()
